# "main correction and comment" - append the next block of sun-track
# readings (16:30 - 20:00) to Sheet1 and leave the selection where data
# entry left off, matching the position check described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New 15-minute-interval rows following row 38 (last existing row, A38=1615).
# Azimut/Altitude stay flat at the tail of the data set (sun has set / track
# clamps), same as B38/C38.
$newRows = @(
    @(1630, 237.6, 2.47),
    @(1645, 237.6, 2.47),
    @(1700, 237.6, 2.47),
    @(1715, 237.6, 2.47),
    @(1730, 237.6, 2.47),
    @(1745, 237.6, 2.47),
    @(1800, 237.6, 2.47),
    @(1815, 237.6, 2.47),
    @(1830, 237.6, 2.47),
    @(1845, 237.6, 2.47),
    @(1900, 237.6, 2.47),
    @(1915, 237.6, 2.47),
    @(1930, 237.6, 2.47),
    @(1945, 237.6, 2.47),
    @(2000, 237.6, 2.47)
)

$startRow = 39
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}

# Scroll/position the view the way it was left after typing the new rows.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("A54").Select() | Out-Null
